$d = $word.ActiveDocument
$table = $d.Tables(1)

function Set-CellText($row, $col, $text) {
    $cell = $table.Cell($row, $col)
    $cell.Range.Text = $text
}

# Row 1 (problems block 1)
Set-CellText 1 1 "24÷5=4, 4"
Set-CellText 1 2 "24÷8=3, 0"
Set-CellText 1 3 "76÷3=25, 1"
Set-CellText 1 4 "55÷3=18, 1"
Set-CellText 1 5 "44÷6=7, 2"

# Row 5 (problems block 2)
Set-CellText 5 1 "38÷4=9, 2"
Set-CellText 5 2 "16÷2=8, 0"
Set-CellText 5 3 "98÷2=49, 0"
Set-CellText 5 4 "94÷3=31, 1"
Set-CellText 5 5 "77÷4=19, 1"

# Row 9 (problems block 3)
Set-CellText 9 1 "34÷2=17, 0"
Set-CellText 9 2 "34÷4=8, 2"
Set-CellText 9 3 "42÷6=7, 0"
Set-CellText 9 4 "72÷3=24, 0"
Set-CellText 9 5 "41÷9=4, 5"

# Row 13 (problems block 4)
Set-CellText 13 1 "20÷6=3, 2"
Set-CellText 13 2 "64÷8=8, 0"
Set-CellText 13 3 "23÷9=2, 5"
Set-CellText 13 4 "61÷4=15, 1"
Set-CellText 13 5 "96÷3=32, 0"

# Row 17 (problems block 5)
Set-CellText 17 1 "84÷6=14, 0"
Set-CellText 17 2 "42÷5=8, 2"
Set-CellText 17 3 "35÷4=8, 3"
Set-CellText 17 4 "98÷4=24, 2"
Set-CellText 17 5 "22÷7=3, 1"
